# 20 setting up run modes
# Adds a new "TestSuite" sheet between addCustomerTest and openAccountTest, appends
# more test-data rows to addCustomerTest, and updates openAccountTest with an extra
# AlertText column + refreshed sample values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. addCustomerTest -- append three more data rows
# ---------------------------------------------------------------------------
$wsAdd = $wb.Worksheets.Item("addCustomerTest")

$wsAdd.Range("A3").Value = "Radha"
$wsAdd.Range("B3").Value = "Madhusudan"
$wsAdd.Range("C3").Value = 354357
$wsAdd.Range("D3").Value = "Customer added successfully"

$wsAdd.Range("A4").Value = "Radha"
$wsAdd.Range("B4").Value = "Rasbihari"
$wsAdd.Range("C4").Value = 354358
$wsAdd.Range("D4").Value = "Customer added successfully"

$wsAdd.Range("A5").Value = "Govinda"
$wsAdd.Range("B5").Value = "Gopal"
$wsAdd.Range("C5").Value = 354359
$wsAdd.Range("D5").Value = "Customer added successfully"

# Wrap the header text + give every populated row the same taller layout.
$wsAdd.Range("A1:D5").WrapText = $true
$wsAdd.Range("A1:D5").RowHeight = 30

$wsAdd.Range("F4").Select()

# ---------------------------------------------------------------------------
# 2. Insert the new "TestSuite" sheet right before openAccountTest
# ---------------------------------------------------------------------------
$wsSuite = $wb.Worksheets.Add($wb.Worksheets.Item("openAccountTest"))
$wsSuite.Name = "TestSuite"

$wsSuite.Range("A1").Value = "TCID"
$wsSuite.Range("B1").Value = "RUNMODE"

$wsSuite.Range("A2").Value = "LoginTest"
$wsSuite.Range("B2").Value = "Y"

$wsSuite.Range("A3").Value = "AddCustomerTest"
$wsSuite.Range("B3").Value = "Y"

$wsSuite.Range("A4").Value = "OpenAccountTest"
$wsSuite.Range("B4").Value = "N"

$wsSuite.Range("A1:B4").WrapText = $true
$wsSuite.Range("A2:B2").RowHeight = 30
$wsSuite.Range("A3:B3").RowHeight = 45
$wsSuite.Range("A4:B4").RowHeight = 30

$wsSuite.Columns.Item(2).ColumnWidth = 9.592447916666666

$wsSuite.Range("B5").Select()
$wsSuite.Activate()

# ---------------------------------------------------------------------------
# 3. openAccountTest -- add AlertText column + refresh sample values
# ---------------------------------------------------------------------------
# Re-resolve the sheet by name: inserting TestSuite shifted the worksheet
# positions, and a stale reference captured before the insert would now
# silently point at the new sheet instead.
$wsOpen = $wb.Worksheets.Item("openAccountTest")

$wsOpen.Range("C1").Value = "AlertText"

$wsOpen.Range("A2").Value = "Radha Raman"
$wsOpen.Range("B2").Value = "Rupee"
$wsOpen.Range("C2").Value = "Account created successfully with account Number "

$wsOpen.Range("A1:C2").WrapText = $true
$wsOpen.Range("A1:C1").RowHeight = 30
$wsOpen.Range("A2:C2").RowHeight = 90

$wsOpen.Columns.Item(1).ColumnWidth = 13.022135416666666

$wsOpen.Range("F7").Select()
